# Applies the "Saldo" export update:
#  - removes the WASHINGTON (005231126) row
#  - adds a BRUNO (004515341) row and a JOSE (004639776) row after RAFAEL (004454365)
#  - adds a LARA (004452597) row after HFR (004361159)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the WASHINGTON row (row 3: 005231126 / WASHINGTON / 14338.69)
$ws.Rows.Item(3).Delete()

# After the delete, row layout is:
#   1 Conta/Nome/Saldo (header)
#   2 004224011 THOMAS     21582.35
#   3 004454365 RAFAEL     13735.23
#   4 004361159 HFR         5714.31
#   5 004488571 CARLOS      1000
#   ...

# 2) Insert BRUNO and JOSE rows right after RAFAEL (row 3), before HFR (row 4)
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "'004515341"
$ws.Range("B4").Value = "BRUNO"
$ws.Range("C4").Value = 13374.08

$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "'004639776"
$ws.Range("B5").Value = "JOSE"
$ws.Range("C5").Value = 12000

# HFR is now row 6. 3) Insert LARA row right after HFR (row 6)
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "'004452597"
$ws.Range("B7").Value = "LARA"
$ws.Range("C7").Value = 1984.11
